$d = $word.ActiveDocument

# --- Paragraph 1: "At country level" / ":" -> add italic (w:i / w:iCs) to
#     paragraph mark rPr and to both runs' rPr.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("At country level", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find paragraph 'At country level'"
}
$p1 = $rng1.Paragraphs(1)
$pr1 = $p1.Range

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' +
    '<w:p w14:paraId="578E5905" w14:textId="01FF90F2" w:rsidR="006A09C1" w:rsidRPr="00F12752" w:rsidRDefault="006A09C1" w:rsidP="00F12752">' +
      '<w:pPr>' +
        '<w:pStyle w:val="MText"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="30"/></w:numPr>' +
        '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr>' +
      '</w:pPr>' +
      '<w:r w:rsidRPr="00F12752">' +
        '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr>' +
        '<w:t>At country level</w:t>' +
      '</w:r>' +
      '<w:r w:rsidR="00F12752">' +
        '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr>' +
        '<w:t>:</w:t>' +
      '</w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$pr1.InsertXML($xml1) | Out-Null

# --- Paragraph 2: "At regional and global levels" -> add italic to
#     paragraph mark rPr and to the run, plus append a new bold+italic
#     run containing ":".
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("At regional and global levels", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find paragraph 'At regional and global levels'"
}
$p2 = $rng2.Paragraphs(1)
$pr2 = $p2.Range

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' +
    '<w:p w14:paraId="48A20F34" w14:textId="77777777" w:rsidR="006A09C1" w:rsidRPr="00F12752" w:rsidRDefault="006A09C1" w:rsidP="00F12752">' +
      '<w:pPr>' +
        '<w:pStyle w:val="MText"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="30"/></w:numPr>' +
        '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr>' +
      '</w:pPr>' +
      '<w:r w:rsidRPr="00F12752">' +
        '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr>' +
        '<w:t>At regional and global levels</w:t>' +
      '</w:r>' +
      '<w:r>' +
        '<w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr>' +
        '<w:t>:</w:t>' +
      '</w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$pr2.InsertXML($xml2) | Out-Null

Write-Output "Applied italic formatting tweaks."
